# Implementing fuel price projections for LCOH
#
# 1) "parameters" sheet: rename the O&M sensitivity columns from
#    high/low to best/worst, and add a new "WI" scenario row.
# 2) "readme" sheet: document the projection year for the fuel price
#    inputs (new "Year" column) and fix the electricity price units
#    label.
#
# NOTE: new literal strings are introduced in the same order the
# original author's edit introduced them, so that the shared-string
# table comes out in the same order as the reference workbook.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("parameters")
$ws2 = $wb.Worksheets.Item("readme")

# --- parameters: rename O&M columns (high/low -> best/worst) --------------
$ws1.Range("F1").Value = "ngboiler_om_best"
$ws1.Range("G1").Value = "ngboiler_om_worst"
$ws1.Range("H1").Value = "eboiler_om_best"
$ws1.Range("I1").Value = "eboiler_om_worst"
$ws1.Range("J1").Value = "hthp_om_best"
$ws1.Range("K1").Value = "hthp_om_worst"

# --- readme: electricity price units, "cents per kwh" -> "$ per kWh" ------
$ws2.Hyperlinks.Delete()
$ws2.Range("C3").Value = "$ per kWh"

# --- parameters: new WI scenario row ---------------------------------------
$ws1.Range("A5").Value = "WI"
$ws1.Range("B5").Value = 0.065
$ws1.Range("C5").Value = 30
$ws1.Range("D5").Value = 0.0852
$ws1.Range("E5").Value = 5.76
$ws1.Range("F5").Value = 0.03
$ws1.Range("G5").Value = 0.06
$ws1.Range("H5").Value = 0.01
$ws1.Range("I5").Value = 0.01
$ws1.Range("J5").Value = 0.01
$ws1.Range("K5").Value = 0.05
$ws1.Range("B5:K5").Font.Color = 0
[void]$ws1.Range("E6").Select()

# --- readme: insert a "Year" column between Units and Source --------------
$ws2.Columns.Item(4).Insert()
$ws2.Columns.Item(4).ColumnWidth = $ws2.Columns.Item(3).ColumnWidth

$ws2.Range("D2").Value = "Year"
$ws2.Range("D3").Value = 2024
$ws2.Range("D4").Value = 2024

# --- readme: hyperlink moved from the old D4 to the new E4 ----------------
$ws2.Hyperlinks.Add($ws2.Range("E4"), "https://www.eia.gov/dnav/ng/ng_pri_sum_a_EPG0_PIN_DMcf_a.htm")
$ws2.Range("E4").Style = "Hyperlink"
